$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.374.01"
$ws.Range("E2").Value = "  +12.39%  "

$ws.Range("D3").Value = "1.816.66"
$ws.Range("E3").Value = "  +7.62%  "

$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.58"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.545"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.24%  "

$ws.Range("E7").Value = "  +0.39%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.35"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.48"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.70%  "

$ws.Range("E10").Value = "  +6.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0679"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0933"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.39%  "

$ws.Range("D13").Value = "2.079.48"
$ws.Range("E13").Value = "  +7.65%  "

$ws.Range("D14").Value = "1.804.63"
$ws.Range("E14").Value = "  +7.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.643"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.77%  "

$ws.Range("D16").Value = "34.358.22"
$ws.Range("E16").Value = "  +12.27%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.34"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.86%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.19"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.27"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "260.19"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.60%  "

$ws.Range("D21").Value = "0.0₃0750"
$ws.Range("E21").Value = "  +4.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.50"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.38"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.41%  "

$ws.Range("E25").Value = "  -0.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.15"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.77"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.116"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.46%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.14"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.30%  "

$ws.Range("E30").Value = "  +0.32%  "

$ws.Range("E31").Value = "  +9.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0515"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.96%  "

$ws.Range("E33").Value = "  +6.85%  "

$ws.Range("E34").Value = "  +7.97%  "

$ws.Range("D35").Value = "1.578.18"
$ws.Range("E35").Value = "  +4.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.83"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.07"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.94%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0189"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.630"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.54%  "

$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "84.91"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.51%  "

$ws.Range("E41").Value = "  +4.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.35"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.915"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.13"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0521"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.44%  "

$ws.Range("E46").Value = "  +4.24%  "

$ws.Range("D47").Value = "1.977.21"
$ws.Range("E47").Value = "  +8.13%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.74"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.45%  "

$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.11"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.18%  "

$ws.Range("E50").Value = "  +0.26%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0122"
$ws.Range("E51").Value = "  +5.44%  "
